$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.28"
$ws.Range("E2").Value = "'-3.10%"
$ws.Range("D3").Value = "'53.53"
$ws.Range("E3").Value = "'7.96%"
$ws.Range("D4").Value = "'5.130"
$ws.Range("E4").Value = "'-2.44%"
$ws.Range("D5").Value = "'0.07842"
$ws.Range("E5").Value = "'-1.48%"
$ws.Range("D6").Value = "'4.525"
$ws.Range("E6").Value = "'-1.06%"
$ws.Range("E7").Value = "'-3.77%"
$ws.Range("D8").Value = "'1.574"
$ws.Range("E8").Value = "'-4.08%"
$ws.Range("D9").Value = "'0.1218"
$ws.Range("E9").Value = "'-6.26%"
$ws.Range("D10").Value = "'0.2008"
$ws.Range("E10").Value = "'2.03%"
$ws.Range("D11").Value = "'0.04724"
$ws.Range("E11").Value = "'2.04%"
$ws.Range("D12").Value = "'0.09446"
$ws.Range("E12").Value = "'-0.64%"
$ws.Range("D13").Value = "'0.1044"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("D14").Value = "'0.001261"
$ws.Range("E14").Value = "'-4.51%"
$ws.Range("D15").Value = "'0.005784"
$ws.Range("E15").Value = "'-1.58%"
$ws.Range("E16").Value = "'2,019.02%"
$ws.Range("D17").Value = "'3.336"
$ws.Range("E17").Value = "'-0.23%"
$ws.Range("D18").Value = "'2.413"
$ws.Range("E18").Value = "'-0.88%"
$ws.Range("D19").Value = "'0.3419"
$ws.Range("E19").Value = "'-0.72%"
$ws.Range("D20").Value = "'8.032"
$ws.Range("E20").Value = "'-2.39%"
$ws.Range("D21").Value = "'0.1370"
$ws.Range("E21").Value = "'-1.41%"
$ws.Range("D22").Value = "'0.3088"
$ws.Range("E22").Value = "'-0.02%"
$ws.Range("D23").Value = "'0.04165"
$ws.Range("E23").Value = "'0.19%"
$ws.Range("D24").Value = "'0.001258"
$ws.Range("E24").Value = "'-4.21%"
$ws.Range("D25").Value = "'0.003913"
$ws.Range("E25").Value = "'-8.05%"
$ws.Range("D26").Value = "'0.0001345"
$ws.Range("E26").Value = "'-0.16%"
$ws.Range("D38").Value = "'0.02604"
$ws.Range("E38").Value = "'-2.76%"
$ws.Range("E39").Value = "'0.98%"
$ws.Range("D40").Value = "'0.01051"
$ws.Range("E40").Value = "'-3.66%"
$ws.Range("D41").Value = "'0.007945"
$ws.Range("E41").Value = "'-0.69%"
$ws.Range("D42").Value = "'0.1423"
$ws.Range("E42").Value = "'-1.13%"
$ws.Range("D43").Value = "'0.008210"
$ws.Range("E43").Value = "'6.17%"
$ws.Range("D44").Value = "'0.008450"
$ws.Range("E44").Value = "'-2.68%"
$ws.Range("D45").Value = "'0.3127"
$ws.Range("E45").Value = "'-2.01%"
$ws.Range("D46").Value = "'0.00007218"
$ws.Range("E46").Value = "'9.00%"
$ws.Range("D47").Value = "'0.00000000747"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("D48").Value = "'0.05330"
$ws.Range("E48").Value = "'-3.09%"
$ws.Range("D49").Value = "'0.002611"
$ws.Range("E49").Value = "'-34.63%"
$ws.Range("D50").Value = "'0.00002093"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0001993"
$ws.Range("E51").Value = "'-0.19%"
